$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (Price / Volume) to Text format so that
# numeric-looking values are stored as literal text, matching the
# original workbook layout (inline strings like "329.94", "2.99%").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "329.94"
$ws.Range("E2").Value = "2.99%"
$ws.Range("D3").Value = "41.13"
$ws.Range("E3").Value = "4.68%"
$ws.Range("D4").Value = "5.681"
$ws.Range("E4").Value = "-3.55%"
$ws.Range("D5").Value = "0.08171"
$ws.Range("E5").Value = "2.01%"
$ws.Range("D6").Value = "2.073"
$ws.Range("E6").Value = "9.38%"
$ws.Range("D7").Value = "8.734"
$ws.Range("E7").Value = "1.20%"
$ws.Range("D8").Value = "4.542"
$ws.Range("E8").Value = "-0.50%"
$ws.Range("E9").Value = "0.42%"
$ws.Range("D10").Value = "0.9260"
$ws.Range("E10").Value = "-0.94%"
$ws.Range("D11").Value = "0.1255"
$ws.Range("E11").Value = "0.42%"
$ws.Range("D12").Value = "0.1960"
$ws.Range("E12").Value = "0.49%"
$ws.Range("D13").Value = "0.09409"
$ws.Range("E13").Value = "3.08%"
$ws.Range("D14").Value = "0.03691"
$ws.Range("E14").Value = "4.99%"
$ws.Range("E15").Value = "10.32%"
$ws.Range("D16").Value = "0.001314"
$ws.Range("E16").Value = "1.66%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.006244"
$ws.Range("E17").Value = "0.39%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.416"
$ws.Range("E18").Value = "1.84%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3485"
$ws.Range("E19").Value = "-1.43%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "8.315"
$ws.Range("E20").Value = "-4.82%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "0.1380"
$ws.Range("E21").Value = "-3.54%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "0.2653"
$ws.Range("E22").Value = "10.00%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "0.04442"
$ws.Range("E23").Value = "-0.45%"
$ws.Range("D24").Value = "0.001270"
$ws.Range("E24").Value = "0.64%"
$ws.Range("D25").Value = "0.004331"
$ws.Range("E25").Value = "-0.88%"
$ws.Range("E26").Value = "3.67%"
$ws.Range("D39").Value = "0.02791"
$ws.Range("E39").Value = "16.59%"
$ws.Range("D40").Value = "0.05482"
$ws.Range("E40").Value = "6.10%"
$ws.Range("D41").Value = "0.007675"
$ws.Range("E41").Value = "2.85%"
$ws.Range("D42").Value = "0.009435"
$ws.Range("E42").Value = "2.95%"
$ws.Range("D43").Value = "0.1418"
$ws.Range("E43").Value = "1.11%"
$ws.Range("D44").Value = "0.002134"
$ws.Range("E44").Value = "0.63%"
$ws.Range("D45").Value = "0.01102"
$ws.Range("E45").Value = "-1.14%"
$ws.Range("D46").Value = "0.00006875"
$ws.Range("E46").Value = "1.85%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.15%"
$ws.Range("D48").Value = "0.002284"
$ws.Range("E48").Value = "60.54%"
$ws.Range("D49").Value = "0.003234"
$ws.Range("E49").Value = "7.58%"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "0.15%"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").Value = "0.15%"

# Remove the temporary text-format marker so the cells end up with
# no explicit style, matching the source workbook.
$ws.Range("D2:E51").ClearFormats()
